# Update "F" column (想去人数 / interested-count) values produced by a
# re-run of the scraper. Sheet names:
#   展览    -> sheet1.xml
#   演出    -> sheet2.xml (unchanged)
#   本地生活 -> sheet3.xml
#   全部类型 -> sheet4.xml

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetLocalLife   = $wb.Worksheets.Item("本地生活")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$sheetExhibition.Range("F6").Value  = 7695
$sheetExhibition.Range("F13").Value = 1709
$sheetExhibition.Range("F15").Value = 6132
$sheetExhibition.Range("F17").Value = 2354
$sheetExhibition.Range("F41").Value = 1164
$sheetExhibition.Range("F42").Value = 473
$sheetExhibition.Range("F44").Value = 3170
$sheetExhibition.Range("F46").Value = 398
$sheetExhibition.Range("F47").Value = 32

# 本地生活 (sheet3)
$sheetLocalLife.Range("F9").Value  = 2099
$sheetLocalLife.Range("F10").Value = 8832
$sheetLocalLife.Range("F11").Value = 948

# 全部类型 (sheet4)
$sheetAllTypes.Range("F5").Value  = 7695
$sheetAllTypes.Range("F8").Value  = 948
$sheetAllTypes.Range("F18").Value = 6132
$sheetAllTypes.Range("F39").Value = 3170
$sheetAllTypes.Range("F43").Value = 32
